# Commit: "add OTree derived  classes"
#
# The underlying data change is on the "datalist" sheet (the workbook's
# ActiveSheet): two of the long fish "info" description strings (electric
# catfish / row 5, and pike / row 16) are shortened by dropping their
# trailing clause, a handful of numeric stat cells are tweaked, and the
# active-cell selection moves from N3 to F13. Re-pointing the two info
# cells at new text automatically retires the old shared-string entries and
# appends the new ones at the end of the table, which is what reshuffles
# every other shared-string index in the diff - so we don't need to touch
# the shared string table by hand, just set the cell values/selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- shortened "info" (column N) descriptions -----------------------------

# Row 5 = electric_catfish: drop "，可做为食用鱼、游钓鱼及观赏鱼" from the end.
$ws.Cells.Item(5, 14).Value = "电鲶，为辐鳍鱼纲鲶形目电鲶科的其中一种，分布于非洲尼罗河、查德湖、尼日河、塞内加尔河、图尔卡纳湖等流域，体长可达122公分，生活在岩石、树根沉积的底中层水域，在夜间活动，属肉食性，具有发电器官，以电击击昏猎物。"

# Row 16 = pike: drop "，适合各种烹饪方式食用" from the end.
$ws.Cells.Item(16, 14).Value = "白斑狗鱼，为辐鳍鱼纲狗鱼目狗鱼科的其中一种。分布于北美洲及欧亚大陆74°N-36°N的淡水流域，体长可达137公分，主要栖息在有植被生长的泠水湖泊、河川，属肉食性，以鱼类、甲壳类、鸟类、小型哺乳类等为食，可作为观赏鱼、游钓鱼及食用鱼。"

# --- numeric stat tweaks ----------------------------------------------------

# Row 3 = carp: maxPopulation (L) 1 -> 3
$ws.Cells.Item(3, 12).Value = 3

# Row 4 = discus: maxPopulation (L) 1 -> 3
$ws.Cells.Item(4, 12).Value = 3

# Row 11 = moorish_idol: speedDiff (E) 6 -> 7
$ws.Cells.Item(11, 5).Value = 7

# Row 14 = oreochromis: depth (C) 1 -> 0, maxPopulation (L) 2 -> 3
$ws.Cells.Item(14, 3).Value = 0
$ws.Cells.Item(14, 12).Value = 3

# Row 16 = pike: avgSpeed (D) 10 -> 11, speedDiff (E) 6 -> 7, maxPopulation (L) 3 -> 2
$ws.Cells.Item(16, 4).Value = 11
$ws.Cells.Item(16, 5).Value = 7
$ws.Cells.Item(16, 12).Value = 2

# --- selection moves from N3 to F13 ----------------------------------------

$ws.Range("F13").Select()
